$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 22:22"

# Row 8 - Alemania
$ws.Range("B8").Value = 150062
$ws.Range("C8").Value = 1609
$ws.Range("E8").Value = 45412
$ws.Range("G8").Value = 164
$ws.Range("H8").Value = 5250

# Row 18 - Suiza
$ws.Range("D18").Value = 19900
$ws.Range("E18").Value = 6859

# Row 88 - Tunez
$ws.Range("B88").Value = 909
$ws.Range("C88").Value = 8
$ws.Range("D88").Value = 190
$ws.Range("E88").Value = 681
$ws.Range("F88").Value = 32

# Row 97 - Niger
$ws.Range("B97").Value = 662
$ws.Range("C97").Value = 5
$ws.Range("D97").Value = 193
$ws.Range("E97").Value = 447
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = 22

$wb.Save()
